# Prepare pseudo random forest experiments
# Update the ODTE sheet's base-estimator hyperparameter example from an
# SVM "rbf" kernel config to a DecisionTree-like "liblinear" kernel with
# a "multiclass_strategy": "ovr" entry instead of "max_features": "sqrt".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ODTE")
$ws.Activate()

# Swap which optional column is populated: multiclass_strategy (K14) now
# holds "ovr" instead of max_features (L14) holding "sqrt".
$ws.Range("K14").Value = """ovr"""
$ws.Range("L14").Value = ""

# C14 holds the "kernel" example value (was "rbf").
$ws.Range("C14").Value = """liblinear"""

# Match the author's updated selection/view state on the sheet.
$ws.Range("E13").Select()
